$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the execution time for "js in browser" (row 2, column C) to "3h 23m"
$ws.Range("C2").Value = "3h 23m"

# Move/update the active selection to C2 (as reflected in the saved sheet view)
$ws.Range("C2").Select()
